$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap row 4 and row 5
$ws.Range("B4").Value = 5474751
$ws.Range("B5").Value = 5475221
$ws.Range("F4").Value = 'AD San Juan'
$ws.Range("F5").Value = 'Racing Rioja CF'
$ws.Range("G4").Value = 'Gernika'
$ws.Range("G5").Value = 'SD Tarazona'
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J4").Value = 'A'
$ws.Range("J5").Value = 'H'
$ws.Range("K4").Value = 2.15
$ws.Range("K5").Value = 2.5
$ws.Range("L4").Value = 3
$ws.Range("L5").Value = 3
$ws.Range("M4").Value = 3.2
$ws.Range("M5").Value = 2.6
$ws.Range("N4").Value = 2.55
$ws.Range("N5").Value = 2.6
$ws.Range("O4").Value = 3
$ws.Range("O5").Value = 3.1
$ws.Range("P4").Value = 2.7
$ws.Range("P5").Value = 2.45
$ws.Range("Q4").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R4").Value = 1.85
$ws.Range("R5").Value = 1.975
$ws.Range("S4").Value = 1.95
$ws.Range("S5").Value = 1.825
$ws.Range("T4").Value = 2
$ws.Range("T5").Value = 2
$ws.Range("U4").Value = 1.875
$ws.Range("U5").Value = 1.8
$ws.Range("V4").Value = 1.925
$ws.Range("V5").Value = 2
$ws.Range("W4").Value = -1
$ws.Range("W5").Value = 1.6
$ws.Range("X4").Value = -1
$ws.Range("X5").Value = -1
$ws.Range("Y4").Value = 1.7
$ws.Range("Y5").Value = -1
$ws.Range("Z4").Value = -1
$ws.Range("Z5").Value = 0.9750000000000001
$ws.Range("AA4").Value = 0.95
$ws.Range("AA5").Value = -1
$ws.Range("AB4").Value = -1
$ws.Range("AB5").Value = 0
$ws.Range("AC4").Value = 0.925
$ws.Range("AC5").Value = -0

# Swap row 7 and row 8
$ws.Range("B7").Value = 5474753
$ws.Range("B8").Value = 5475220
$ws.Range("F7").Value = 'Arnedo'
$ws.Range("F8").Value = 'Mutilvera'
$ws.Range("G7").Value = 'Utebo'
$ws.Range("G8").Value = 'Arenas Club de Getxo'
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("J7").Value = 'D'
$ws.Range("J8").Value = 'H'
$ws.Range("K7").Value = 3.6
$ws.Range("K8").Value = 2.5
$ws.Range("L7").Value = 3.2
$ws.Range("L8").Value = 2.8
$ws.Range("M7").Value = 1.909
$ws.Range("M8").Value = 2.8
$ws.Range("N7").Value = 3.8
$ws.Range("N8").Value = 2.375
$ws.Range("O7").Value = 3.3
$ws.Range("O8").Value = 2.75
$ws.Range("P7").Value = 1.833
$ws.Range("P8").Value = 3.1
$ws.Range("Q7").Value = 0.5
$ws.Range("Q8").Value = -0.25
$ws.Range("R7").Value = 1.85
$ws.Range("R8").Value = 2.05
$ws.Range("S7").Value = 1.95
$ws.Range("S8").Value = 1.75
$ws.Range("T7").Value = 2.25
$ws.Range("T8").Value = 2
$ws.Range("U7").Value = 1.9
$ws.Range("U8").Value = 1.825
$ws.Range("V7").Value = 1.9
$ws.Range("V8").Value = 1.975
$ws.Range("W7").Value = -1
$ws.Range("W8").Value = 1.375
$ws.Range("X7").Value = 2.3
$ws.Range("X8").Value = -1
$ws.Range("Y7").Value = -1
$ws.Range("Y8").Value = -1
$ws.Range("Z7").Value = 0.8500000000000001
$ws.Range("Z8").Value = 1.05
$ws.Range("AA7").Value = -1
$ws.Range("AA8").Value = -1
$ws.Range("AB7").Value = -0.5
$ws.Range("AB8").Value = -1
$ws.Range("AC7").Value = 0.45
$ws.Range("AC8").Value = 0.9750000000000001

# Swap row 18 and row 19
$ws.Range("B18").Value = 5466417
$ws.Range("B19").Value = 5474760
$ws.Range("F18").Value = 'Tudelano'
$ws.Range("F19").Value = 'SD Tarazona'
$ws.Range("G18").Value = 'SD Beasain'
$ws.Range("G19").Value = 'Real Sociedad C'
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 2
$ws.Range("I18").Value = 2
$ws.Range("I19").Value = 1
$ws.Range("J18").Value = 'A'
$ws.Range("J19").Value = 'H'
$ws.Range("K18").Value = 1.55
$ws.Range("K19").Value = 1.909
$ws.Range("L18").Value = 3.75
$ws.Range("L19").Value = 3
$ws.Range("M18").Value = 5
$ws.Range("M19").Value = 3.9
$ws.Range("N18").Value = 1.727
$ws.Range("N19").Value = 1.8
$ws.Range("O18").Value = 3.5
$ws.Range("O19").Value = 3.2
$ws.Range("P18").Value = 3.8
$ws.Range("P19").Value = 4.2
$ws.Range("Q18").Value = -0.75
$ws.Range("Q19").Value = -0.5
$ws.Range("R18").Value = 2
$ws.Range("R19").Value = 1.875
$ws.Range("S18").Value = 1.8
$ws.Range("S19").Value = 1.925
$ws.Range("T18").Value = 2
$ws.Range("T19").Value = 2.25
$ws.Range("U18").Value = 1.9
$ws.Range("U19").Value = 1.85
$ws.Range("V18").Value = 1.9
$ws.Range("V19").Value = 1.95
$ws.Range("W18").Value = -1
$ws.Range("W19").Value = 0.8
$ws.Range("X18").Value = -1
$ws.Range("X19").Value = -1
$ws.Range("Y18").Value = 2.8
$ws.Range("Y19").Value = -1
$ws.Range("Z18").Value = -1
$ws.Range("Z19").Value = 0.875
$ws.Range("AA18").Value = 0.8
$ws.Range("AA19").Value = -1
$ws.Range("AB18").Value = 0.8999999999999999
$ws.Range("AB19").Value = 0.8500000000000001
$ws.Range("AC18").Value = -1
$ws.Range("AC19").Value = -1

# Swap row 40 and row 41
$ws.Range("B40").Value = 5474769
$ws.Range("B41").Value = 5466420
$ws.Range("F40").Value = 'UD Logrones B'
$ws.Range("F41").Value = 'Utebo'
$ws.Range("G40").Value = 'Real Sociedad C'
$ws.Range("G41").Value = 'Tudelano'
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 1
$ws.Range("I40").Value = 1
$ws.Range("I41").Value = 1
$ws.Range("J40").Value = 'A'
$ws.Range("J41").Value = 'D'
$ws.Range("K40").Value = 2.25
$ws.Range("K41").Value = 1.95
$ws.Range("L40").Value = 3.25
$ws.Range("L41").Value = 3
$ws.Range("M40").Value = 2.8
$ws.Range("M41").Value = 3.75
$ws.Range("N40").Value = 2.3
$ws.Range("N41").Value = 2.15
$ws.Range("O40").Value = 3
$ws.Range("O41").Value = 2.75
$ws.Range("P40").Value = 2.9
$ws.Range("P41").Value = 3.5
$ws.Range("Q40").Value = -0.25
$ws.Range("Q41").Value = -0.25
$ws.Range("R40").Value = 2.025
$ws.Range("R41").Value = 1.875
$ws.Range("S40").Value = 1.775
$ws.Range("S41").Value = 1.925
$ws.Range("T40").Value = 2.25
$ws.Range("T41").Value = 2
$ws.Range("U40").Value = 1.95
$ws.Range("U41").Value = 1.975
$ws.Range("V40").Value = 1.85
$ws.Range("V41").Value = 1.825
$ws.Range("W40").Value = -1
$ws.Range("W41").Value = -1
$ws.Range("X40").Value = -1
$ws.Range("X41").Value = 1.75
$ws.Range("Y40").Value = 1.9
$ws.Range("Y41").Value = -1
$ws.Range("Z40").Value = -1
$ws.Range("Z41").Value = -0.5
$ws.Range("AA40").Value = 0.7749999999999999
$ws.Range("AA41").Value = 0.4625
$ws.Range("AB40").Value = -1
$ws.Range("AB41").Value = 0
$ws.Range("AC40").Value = 0.8500000000000001
$ws.Range("AC41").Value = -0

# Swap row 76 and row 77
$ws.Range("B76").Value = 5451388
$ws.Range("B77").Value = 5475227
$ws.Range("F76").Value = 'Alaves B'
$ws.Range("F77").Value = 'UD Logrones B'
$ws.Range("G76").Value = 'Sestao River'
$ws.Range("G77").Value = 'SD Tarazona'
$ws.Range("H76").Value = 0
$ws.Range("H77").Value = 1
$ws.Range("I76").Value = 0
$ws.Range("I77").Value = 1
$ws.Range("J76").Value = 'D'
$ws.Range("J77").Value = 'D'
$ws.Range("K76").Value = 3
$ws.Range("K77").Value = 3
$ws.Range("L76").Value = 2.8
$ws.Range("L77").Value = 3
$ws.Range("M76").Value = 2.4
$ws.Range("M77").Value = 2.25
$ws.Range("N76").Value = 3.1
$ws.Range("N77").Value = 3.1
$ws.Range("O76").Value = 2.625
$ws.Range("O77").Value = 3.1
$ws.Range("P76").Value = 2.4
$ws.Range("P77").Value = 2.2
$ws.Range("Q76").Value = 0.25
$ws.Range("Q77").Value = 0.25
$ws.Range("R76").Value = 1.675
$ws.Range("R77").Value = 1.85
$ws.Range("S76").Value = 2.05
$ws.Range("S77").Value = 1.95
$ws.Range("T76").Value = 1.75
$ws.Range("T77").Value = 2.25
$ws.Range("U76").Value = 1.825
$ws.Range("U77").Value = 1.9
$ws.Range("V76").Value = 1.975
$ws.Range("V77").Value = 1.9
$ws.Range("W76").Value = -1
$ws.Range("W77").Value = -1
$ws.Range("X76").Value = 1.625
$ws.Range("X77").Value = 2.1
$ws.Range("Y76").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z76").Value = 0.3375
$ws.Range("Z77").Value = 0.425
$ws.Range("AA76").Value = -0.5
$ws.Range("AA77").Value = -0.5
$ws.Range("AB76").Value = -1
$ws.Range("AB77").Value = -0.5
$ws.Range("AC76").Value = 0.9750000000000001
$ws.Range("AC77").Value = 0.45

# Swap row 79 and row 80
$ws.Range("B79").Value = 5474526
$ws.Range("B80").Value = 5474527
$ws.Range("F79").Value = 'Cirbonero'
$ws.Range("F80").Value = 'Brea'
$ws.Range("G79").Value = 'Arnedo'
$ws.Range("G80").Value = 'SD Beasain'
$ws.Range("H79").Value = 3
$ws.Range("H80").Value = 0
$ws.Range("I79").Value = 2
$ws.Range("I80").Value = 0
$ws.Range("J79").Value = 'H'
$ws.Range("J80").Value = 'D'
$ws.Range("K79").Value = 1.666
$ws.Range("K80").Value = 2.1
$ws.Range("L79").Value = 3.25
$ws.Range("L80").Value = 3
$ws.Range("M79").Value = 5
$ws.Range("M80").Value = 3.3
$ws.Range("N79").Value = 1.4
$ws.Range("N80").Value = 2.25
$ws.Range("O79").Value = 3.8
$ws.Range("O80").Value = 2.875
$ws.Range("P79").Value = 7.5
$ws.Range("P80").Value = 3.1
$ws.Range("Q79").Value = -1.25
$ws.Range("Q80").Value = -0.25
$ws.Range("R79").Value = 2.025
$ws.Range("R80").Value = 1.975
$ws.Range("S79").Value = 1.775
$ws.Range("S80").Value = 1.825
$ws.Range("T79").Value = 2.25
$ws.Range("T80").Value = 2
$ws.Range("U79").Value = 1.85
$ws.Range("U80").Value = 1.9
$ws.Range("V79").Value = 1.95
$ws.Range("V80").Value = 1.9
$ws.Range("W79").Value = 0.3999999999999999
$ws.Range("W80").Value = -1
$ws.Range("X79").Value = -1
$ws.Range("X80").Value = 1.875
$ws.Range("Y79").Value = -1
$ws.Range("Y80").Value = -1
$ws.Range("Z79").Value = -0.5
$ws.Range("Z80").Value = -0.5
$ws.Range("AA79").Value = 0.3875
$ws.Range("AA80").Value = 0.4125
$ws.Range("AB79").Value = 0.8500000000000001
$ws.Range("AB80").Value = -1
$ws.Range("AC79").Value = -1
$ws.Range("AC80").Value = 0.8999999999999999

# Swap row 88 and row 89
$ws.Range("B88").Value = 5474791
$ws.Range("B89").Value = 5474790
$ws.Range("F88").Value = 'Izarra'
$ws.Range("F89").Value = 'Utebo'
$ws.Range("G88").Value = 'Gernika'
$ws.Range("G89").Value = 'UD Logrones B'
$ws.Range("H88").Value = 1
$ws.Range("H89").Value = 0
$ws.Range("I88").Value = 2
$ws.Range("I89").Value = 0
$ws.Range("J88").Value = 'A'
$ws.Range("J89").Value = 'D'
$ws.Range("K88").Value = 2.3
$ws.Range("K89").Value = 1.533
$ws.Range("L88").Value = 2.8
$ws.Range("L89").Value = 3.75
$ws.Range("M88").Value = 3.25
$ws.Range("M89").Value = 5
$ws.Range("N88").Value = 2.55
$ws.Range("N89").Value = 1.615
$ws.Range("O88").Value = 2.625
$ws.Range("O89").Value = 3.8
$ws.Range("P88").Value = 3
$ws.Range("P89").Value = 4.5
$ws.Range("Q88").Value = 0
$ws.Range("Q89").Value = -1
$ws.Range("R88").Value = 1.75
$ws.Range("R89").Value = 2.05
$ws.Range("S88").Value = 2.05
$ws.Range("S89").Value = 1.75
$ws.Range("T88").Value = 1.75
$ws.Range("T89").Value = 2.25
$ws.Range("U88").Value = 1.8
$ws.Range("U89").Value = 1.8
$ws.Range("V88").Value = 2
$ws.Range("V89").Value = 2
$ws.Range("W88").Value = -1
$ws.Range("W89").Value = -1
$ws.Range("X88").Value = -1
$ws.Range("X89").Value = 2.8
$ws.Range("Y88").Value = 2
$ws.Range("Y89").Value = -1
$ws.Range("Z88").Value = -1
$ws.Range("Z89").Value = -1
$ws.Range("AA88").Value = 1.05
$ws.Range("AA89").Value = 0.75
$ws.Range("AB88").Value = 0.8
$ws.Range("AB89").Value = -1
$ws.Range("AC88").Value = -1
$ws.Range("AC89").Value = 1

# Swap row 98 and row 100
$ws.Range("B98").Value = 5474530
$ws.Range("B100").Value = 5474796
$ws.Range("F98").Value = 'Cirbonero'
$ws.Range("F100").Value = 'Gernika'
$ws.Range("G98").Value = 'Real Sociedad C'
$ws.Range("G100").Value = 'SD Tarazona'
$ws.Range("H98").Value = 1
$ws.Range("H100").Value = 1
$ws.Range("I98").Value = 2
$ws.Range("I100").Value = 1
$ws.Range("J98").Value = 'A'
$ws.Range("J100").Value = 'D'
$ws.Range("K98").Value = 2.5
$ws.Range("K100").Value = 2.2
$ws.Range("L98").Value = 3
$ws.Range("L100").Value = 3
$ws.Range("M98").Value = 2.625
$ws.Range("M100").Value = 3.1
$ws.Range("N98").Value = 2.2
$ws.Range("N100").Value = 2.15
$ws.Range("O98").Value = 3.1
$ws.Range("O100").Value = 3.1
$ws.Range("P98").Value = 3.1
$ws.Range("P100").Value = 3.1
$ws.Range("Q98").Value = -0.25
$ws.Range("Q100").Value = -0.25
$ws.Range("R98").Value = 1.95
$ws.Range("R100").Value = 1.925
$ws.Range("S98").Value = 1.85
$ws.Range("S100").Value = 1.875
$ws.Range("T98").Value = 2
$ws.Range("T100").Value = 2
$ws.Range("U98").Value = 1.975
$ws.Range("U100").Value = 1.875
$ws.Range("V98").Value = 1.825
$ws.Range("V100").Value = 1.925
$ws.Range("W98").Value = -1
$ws.Range("W100").Value = -1
$ws.Range("X98").Value = -1
$ws.Range("X100").Value = 2.1
$ws.Range("Y98").Value = 2.1
$ws.Range("Y100").Value = -1
$ws.Range("Z98").Value = -1
$ws.Range("Z100").Value = -0.5
$ws.Range("AA98").Value = 0.8500000000000001
$ws.Range("AA100").Value = 0.4375
$ws.Range("AB98").Value = 0.9750000000000001
$ws.Range("AB100").Value = 0
$ws.Range("AC98").Value = -1
$ws.Range("AC100").Value = -0

# Swap row 196 and row 197
$ws.Range("B196").Value = 6973206
$ws.Range("B197").Value = 6972693
$ws.Range("F196").Value = 'Tudelano'
$ws.Range("F197").Value = 'Gernika'
$ws.Range("G196").Value = 'AD San Juan'
$ws.Range("G197").Value = 'Alaves B'
$ws.Range("H196").Value = 3
$ws.Range("H197").Value = 3
$ws.Range("I196").Value = 0
$ws.Range("I197").Value = 1
$ws.Range("J196").Value = 'H'
$ws.Range("J197").Value = 'H'
$ws.Range("K196").Value = 2.375
$ws.Range("K197").Value = 2.8
$ws.Range("L196").Value = 2.8
$ws.Range("L197").Value = 2.8
$ws.Range("M196").Value = 3
$ws.Range("M197").Value = 2.5
$ws.Range("N196").Value = 1.615
$ws.Range("N197").Value = 2.5
$ws.Range("O196").Value = 3.1
$ws.Range("O197").Value = 2.8
$ws.Range("P196").Value = 5.5
$ws.Range("P197").Value = 2.75
$ws.Range("Q196").Value = -0.75
$ws.Range("Q197").Value = 0
$ws.Range("R196").Value = 1.925
$ws.Range("R197").Value = 1.8
$ws.Range("S196").Value = 1.875
$ws.Range("S197").Value = 2
$ws.Range("T196").Value = 2
$ws.Range("T197").Value = 2
$ws.Range("U196").Value = 1.9
$ws.Range("U197").Value = 1.925
$ws.Range("V196").Value = 1.9
$ws.Range("V197").Value = 1.875
$ws.Range("W196").Value = 0.615
$ws.Range("W197").Value = 1.5
$ws.Range("X196").Value = -1
$ws.Range("X197").Value = -1
$ws.Range("Y196").Value = -1
$ws.Range("Y197").Value = -1
$ws.Range("Z196").Value = 0.925
$ws.Range("Z197").Value = 0.8
$ws.Range("AA196").Value = -1
$ws.Range("AA197").Value = -1
$ws.Range("AB196").Value = 0.8999999999999999
$ws.Range("AB197").Value = 0.925
$ws.Range("AC196").Value = -1
$ws.Range("AC197").Value = -1

# Swap row 227 and row 228
$ws.Range("B227").Value = 6972736
$ws.Range("B228").Value = 6972732
$ws.Range("F227").Value = 'Gernika'
$ws.Range("F228").Value = 'Izarra'
$ws.Range("G227").Value = 'Arenas Club de Getxo'
$ws.Range("G228").Value = 'Valle Egues'
$ws.Range("H227").Value = 2
$ws.Range("H228").Value = 0
$ws.Range("I227").Value = 2
$ws.Range("I228").Value = 2
$ws.Range("J227").Value = 'D'
$ws.Range("J228").Value = 'A'
$ws.Range("K227").Value = 2
$ws.Range("K228").Value = 2.15
$ws.Range("L227").Value = 2.8
$ws.Range("L228").Value = 3.1
$ws.Range("M227").Value = 3.9
$ws.Range("M228").Value = 3.1
$ws.Range("N227").Value = 2.25
$ws.Range("N228").Value = 1.85
$ws.Range("O227").Value = 2.75
$ws.Range("O228").Value = 3.1
$ws.Range("P227").Value = 3.4
$ws.Range("P228").Value = 4.333
$ws.Range("Q227").Value = -0.25
$ws.Range("Q228").Value = -0.5
$ws.Range("R227").Value = 1.925
$ws.Range("R228").Value = 1.9
$ws.Range("S227").Value = 1.875
$ws.Range("S228").Value = 1.9
$ws.Range("T227").Value = 2
$ws.Range("T228").Value = 2
$ws.Range("U227").Value = 1.975
$ws.Range("U228").Value = 1.875
$ws.Range("V227").Value = 1.825
$ws.Range("V228").Value = 1.925
$ws.Range("W227").Value = -1
$ws.Range("W228").Value = -1
$ws.Range("X227").Value = 1.75
$ws.Range("X228").Value = -1
$ws.Range("Y227").Value = -1
$ws.Range("Y228").Value = 3.333
$ws.Range("Z227").Value = -0.5
$ws.Range("Z228").Value = -1
$ws.Range("AA227").Value = 0.4375
$ws.Range("AA228").Value = 0.8999999999999999
$ws.Range("AB227").Value = 0.9750000000000001
$ws.Range("AB228").Value = 0
$ws.Range("AC227").Value = -1
$ws.Range("AC228").Value = -0

# Swap row 256 and row 257
$ws.Range("B256").Value = 6972754
$ws.Range("B257").Value = 6962803
$ws.Range("F256").Value = 'AD San Juan'
$ws.Range("F257").Value = 'Barakaldo'
$ws.Range("G256").Value = 'UD Barbastro'
$ws.Range("G257").Value = 'UD Logrones'
$ws.Range("H256").Value = 3
$ws.Range("H257").Value = 1
$ws.Range("I256").Value = 1
$ws.Range("I257").Value = 1
$ws.Range("J256").Value = 'H'
$ws.Range("J257").Value = 'D'
$ws.Range("K256").Value = 2.3
$ws.Range("K257").Value = 2.5
$ws.Range("L256").Value = 2.9
$ws.Range("L257").Value = 2.625
$ws.Range("M256").Value = 3
$ws.Range("M257").Value = 3
$ws.Range("N256").Value = 2.9
$ws.Range("N257").Value = 3.2
$ws.Range("O256").Value = 2.625
$ws.Range("O257").Value = 2.8
$ws.Range("P256").Value = 2.55
$ws.Range("P257").Value = 2.3
$ws.Range("Q256").Value = 0
$ws.Range("Q257").Value = 0
$ws.Range("R256").Value = 2.075
$ws.Range("R257").Value = 2.05
$ws.Range("S256").Value = 1.725
$ws.Range("S257").Value = 1.75
$ws.Range("T256").Value = 2
$ws.Range("T257").Value = 2.25
$ws.Range("U256").Value = 2
$ws.Range("U257").Value = 1.975
$ws.Range("V256").Value = 1.8
$ws.Range("V257").Value = 1.825
$ws.Range("W256").Value = 1.9
$ws.Range("W257").Value = -1
$ws.Range("X256").Value = -1
$ws.Range("X257").Value = 1.8
$ws.Range("Y256").Value = -1
$ws.Range("Y257").Value = -1
$ws.Range("Z256").Value = 1.075
$ws.Range("Z257").Value = 0
$ws.Range("AA256").Value = -1
$ws.Range("AA257").Value = -0
$ws.Range("AB256").Value = 1
$ws.Range("AB257").Value = -0.5
$ws.Range("AC256").Value = -1
$ws.Range("AC257").Value = 0.4125

# Swap row 295 and row 296
$ws.Range("B295").Value = 6972780
$ws.Range("B296").Value = 6972782
$ws.Range("F295").Value = 'Mutilvera'
$ws.Range("F296").Value = 'Izarra'
$ws.Range("G295").Value = 'AD San Juan'
$ws.Range("G296").Value = 'Tudelano'
$ws.Range("H295").Value = 1
$ws.Range("H296").Value = 1
$ws.Range("I295").Value = 0
$ws.Range("I296").Value = 0
$ws.Range("J295").Value = 'H'
$ws.Range("J296").Value = 'H'
$ws.Range("K295").Value = 2.1
$ws.Range("K296").Value = 3.6
$ws.Range("L295").Value = 2.9
$ws.Range("L296").Value = 3
$ws.Range("M295").Value = 3.4
$ws.Range("M296").Value = 2
$ws.Range("N295").Value = 2.45
$ws.Range("N296").Value = 2.8
$ws.Range("O295").Value = 2.875
$ws.Range("O296").Value = 3.1
$ws.Range("P295").Value = 2.8
$ws.Range("P296").Value = 2.375
$ws.Range("Q295").Value = 0
$ws.Range("Q296").Value = 0
$ws.Range("R295").Value = 1.725
$ws.Range("R296").Value = 2.05
$ws.Range("S295").Value = 2.075
$ws.Range("S296").Value = 1.75
$ws.Range("T295").Value = 2
$ws.Range("T296").Value = 2
$ws.Range("U295").Value = 1.9
$ws.Range("U296").Value = 1.8
$ws.Range("V295").Value = 1.9
$ws.Range("V296").Value = 2
$ws.Range("W295").Value = 1.45
$ws.Range("W296").Value = 1.8
$ws.Range("X295").Value = -1
$ws.Range("X296").Value = -1
$ws.Range("Y295").Value = -1
$ws.Range("Y296").Value = -1
$ws.Range("Z295").Value = 0.7250000000000001
$ws.Range("Z296").Value = 1.05
$ws.Range("AA295").Value = -1
$ws.Range("AA296").Value = -1
$ws.Range("AB295").Value = -1
$ws.Range("AB296").Value = -1
$ws.Range("AC295").Value = 0.8999999999999999
$ws.Range("AC296").Value = 1

# Swap row 342 and row 343
$ws.Range("B342").Value = 6972816
$ws.Range("B343").Value = 6973239
$ws.Range("F342").Value = 'Arenas Club de Getxo'
$ws.Range("F343").Value = 'AD San Juan'
$ws.Range("G342").Value = 'CD Calahorra'
$ws.Range("G343").Value = 'Tudelano'
$ws.Range("H342").Value = 0
$ws.Range("H343").Value = 1
$ws.Range("I342").Value = 1
$ws.Range("I343").Value = 1
$ws.Range("J342").Value = 'A'
$ws.Range("J343").Value = 'D'
$ws.Range("K342").Value = 2.1
$ws.Range("K343").Value = 3.6
$ws.Range("L342").Value = 2.9
$ws.Range("L343").Value = 3
$ws.Range("M342").Value = 3.4
$ws.Range("M343").Value = 2
$ws.Range("N342").Value = 2.375
$ws.Range("N343").Value = 3.5
$ws.Range("O342").Value = 2.625
$ws.Range("O343").Value = 3.1
$ws.Range("P342").Value = 3.2
$ws.Range("P343").Value = 2
$ws.Range("Q342").Value = -0.25
$ws.Range("Q343").Value = 0.5
$ws.Range("R342").Value = 2.025
$ws.Range("R343").Value = 1.725
$ws.Range("S342").Value = 1.775
$ws.Range("S343").Value = 1.975
$ws.Range("T342").Value = 1.75
$ws.Range("T343").Value = 2
$ws.Range("U342").Value = 1.975
$ws.Range("U343").Value = 1.8
$ws.Range("V342").Value = 1.825
$ws.Range("V343").Value = 2
$ws.Range("W342").Value = -1
$ws.Range("W343").Value = -1
$ws.Range("X342").Value = -1
$ws.Range("X343").Value = 2.1
$ws.Range("Y342").Value = 2.2
$ws.Range("Y343").Value = -1
$ws.Range("Z342").Value = -1
$ws.Range("Z343").Value = 0.7250000000000001
$ws.Range("AA342").Value = 0.7749999999999999
$ws.Range("AA343").Value = -1
$ws.Range("AB342").Value = -1
$ws.Range("AB343").Value = 0
$ws.Range("AC342").Value = 0.825
$ws.Range("AC343").Value = -0

# Swap row 374 and row 375
$ws.Range("B374").Value = 6972829
$ws.Range("B375").Value = 6972831
$ws.Range("F374").Value = 'Utebo'
$ws.Range("F375").Value = 'UD Logrones'
$ws.Range("G374").Value = 'Izarra'
$ws.Range("G375").Value = 'Real Zaragoza B'
$ws.Range("K374").Value = 1.533
$ws.Range("K375").Value = 1.571
$ws.Range("L374").Value = 3.4
$ws.Range("L375").Value = 3.6
$ws.Range("M374").Value = 6
$ws.Range("M375").Value = 5
$ws.Range("N374").Value = 1.909
$ws.Range("N375").Value = 1.4
$ws.Range("O374").Value = 3
$ws.Range("O375").Value = 3.8
$ws.Range("P374").Value = 4
$ws.Range("P375").Value = 7
$ws.Range("Q374").Value = -0.5
$ws.Range("Q375").Value = -1.25
$ws.Range("R374").Value = 2
$ws.Range("R375").Value = 2
$ws.Range("S374").Value = 1.8
$ws.Range("S375").Value = 1.8
$ws.Range("T374").Value = 2
$ws.Range("T375").Value = 2.25
$ws.Range("U374").Value = 1.95
$ws.Range("U375").Value = 1.975
$ws.Range("V374").Value = 1.85
$ws.Range("V375").Value = 1.825
$ws.Range("W374").Value = 0
$ws.Range("W375").Value = 0
$ws.Range("X374").Value = 0
$ws.Range("X375").Value = 0
$ws.Range("Y374").Value = 0
$ws.Range("Y375").Value = 0
$ws.Range("Z374").Value = 0
$ws.Range("Z375").Value = 0
$ws.Range("AA374").Value = 0
$ws.Range("AA375").Value = 0

# Rotate rows 198 <- 199 <- 200 <- 198
$ws.Range("B198").Value = 6973211
$ws.Range("B199").Value = 6972702
$ws.Range("B200").Value = 6972711
$ws.Range("F198").Value = 'AD San Juan'
$ws.Range("F199").Value = 'Arenas Club de Getxo'
$ws.Range("F200").Value = 'Mutilvera'
$ws.Range("G198").Value = 'Naxara'
$ws.Range("G199").Value = 'Tudelano'
$ws.Range("G200").Value = 'Real Sociedad C'
$ws.Range("H198").Value = 2
$ws.Range("H199").Value = 1
$ws.Range("H200").Value = 0
$ws.Range("I198").Value = 1
$ws.Range("I199").Value = 2
$ws.Range("I200").Value = 0
$ws.Range("J198").Value = 'H'
$ws.Range("J199").Value = 'A'
$ws.Range("J200").Value = 'D'
$ws.Range("K198").Value = 2.4
$ws.Range("K199").Value = 2.4
$ws.Range("K200").Value = 2.15
$ws.Range("L198").Value = 3
$ws.Range("L199").Value = 3
$ws.Range("L200").Value = 3
$ws.Range("M198").Value = 2.75
$ws.Range("M199").Value = 2.75
$ws.Range("M200").Value = 3.1
$ws.Range("N198").Value = 2
$ws.Range("N199").Value = 2.3
$ws.Range("N200").Value = 2.2
$ws.Range("O198").Value = 3.3
$ws.Range("O199").Value = 3
$ws.Range("O200").Value = 3
$ws.Range("P198").Value = 3.3
$ws.Range("P199").Value = 2.9
$ws.Range("P200").Value = 3
$ws.Range("Q198").Value = -0.25
$ws.Range("Q199").Value = -0.25
$ws.Range("Q200").Value = -0.25
$ws.Range("R198").Value = 1.775
$ws.Range("R199").Value = 2.025
$ws.Range("R200").Value = 2
$ws.Range("S198").Value = 2.025
$ws.Range("S199").Value = 1.775
$ws.Range("S200").Value = 1.8
$ws.Range("T198").Value = 2
$ws.Range("T199").Value = 2
$ws.Range("T200").Value = 2.25
$ws.Range("U198").Value = 1.725
$ws.Range("U199").Value = 1.975
$ws.Range("U200").Value = 1.95
$ws.Range("V198").Value = 1.975
$ws.Range("V199").Value = 1.725
$ws.Range("V200").Value = 1.85
$ws.Range("W198").Value = 1
$ws.Range("W199").Value = -1
$ws.Range("W200").Value = -1
$ws.Range("X198").Value = -1
$ws.Range("X199").Value = -1
$ws.Range("X200").Value = 2
$ws.Range("Y198").Value = -1
$ws.Range("Y199").Value = 1.9
$ws.Range("Y200").Value = -1
$ws.Range("Z198").Value = 0.7749999999999999
$ws.Range("Z199").Value = -1
$ws.Range("Z200").Value = -0.5
$ws.Range("AA198").Value = -1
$ws.Range("AA199").Value = 0.7749999999999999
$ws.Range("AA200").Value = 0.4
$ws.Range("AB198").Value = 0.7250000000000001
$ws.Range("AB199").Value = 0.9750000000000001
$ws.Range("AB200").Value = -1
$ws.Range("AC198").Value = -1
$ws.Range("AC199").Value = -1
$ws.Range("AC200").Value = 0.8500000000000001
